$d = $word.ActiveDocument

# Locate the single occurrence of "facebook" inside the OBSERVATION
# paragraph ("... is clicked it should display the facebook page of
# Naveen Engineering &Speciality Coating Page should be displayed").
$target = $d.Content.Duplicate
$target.Find.Execute("facebook", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)

$start = $target.Start
$end = $target.End

# Word normally re-merges a freshly retyped span back into its
# neighbouring runs when the resulting formatting is identical. Drop a
# temporary bookmark at each split point so the boundaries survive the
# edit, producing three separate (but identically formatted) runs:
#   " is clicked it should display the "
#   "respective social media"
#   " page of Naveen Engineering &Speciality Coating Page should be displayed"
$d.Bookmarks.Add("zzSplitStart", $d.Range($start, $start))
$d.Bookmarks.Add("zzSplitEnd", $d.Range($end, $end))

$mid = $d.Range($start, $end)
$mid.Text = "respective social media"

$d.Bookmarks("zzSplitStart").Delete()
$d.Bookmarks("zzSplitEnd").Delete()
